$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "PUBLONS009" worksheet right after "Test Cases" (i.e.
#    before the existing "PUBLONS005" sheet). Worksheets.Add($null, <after>)
#    places the new sheet immediately after the given sheet.
# ---------------------------------------------------------------------------
$testCases = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $testCases)
$newSheet.Name = "PUBLONS009"

# Re-fetch sheet references by name AFTER the insert - cached references
# taken before a sheet-count-changing operation can resolve to the wrong
# sheet afterwards.
$publons005 = $wb.Worksheets.Item("PUBLONS005")

# ---------------------------------------------------------------------------
# 2. Populate the new sheet's data (values first, then copy formats only
#    from equivalently-styled cells on PUBLONS005 so we land on the exact
#    same style indices instead of synthesizing brand-new ones).
# ---------------------------------------------------------------------------
$newSheet.Range("A1").Value2 = "CHARACTER LENGTH"
$newSheet.Range("B1").Value2 = "VALIDITY"
$newSheet.Range("C1").Value2 = "Runmode"
$newSheet.Range("D1").Value2 = "PASS"

$newSheet.Range("A2").Value2 = 49
$newSheet.Range("B2").Value2 = "YES"
$newSheet.Range("C2").Value2 = "Y"
$newSheet.Range("D2").Value2 = "SKIP"

$newSheet.Range("A3").Value2 = 50
$newSheet.Range("B3").Value2 = "YES"
$newSheet.Range("C3").Value2 = "Y"
$newSheet.Range("D3").Value2 = "SKIP"

$newSheet.Range("A4").Value2 = 51
$newSheet.Range("B4").Value2 = "NO"
$newSheet.Range("C4").Value2 = "Y"
$newSheet.Range("D4").Value2 = "PASS"

# Row 1 header formatting (style index 10 on PUBLONS005!A1:F1).
$publons005.Range("A1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Column B (rows 2-4) uses style 11 (PUBLONS005!D2, also style 11).
$publons005.Range("D2").Copy()
$newSheet.Range("B2:B4").PasteSpecial(-4122)

# Columns A, C, D (rows 2-4) use the plain bordered style 2
# (PUBLONS005!E2 / F2, both style 2).
$publons005.Range("E2").Copy()
$newSheet.Range("A2:A4").PasteSpecial(-4122)
$newSheet.Range("C2:C4").PasteSpecial(-4122)
$publons005.Range("F2").Copy()
$newSheet.Range("D2:D4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Match the column widths / view state captured in the target workbook.
$newSheet.Range("F15").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Update the workbook.xml sheet order bookkeeping is handled automatically
#    by Worksheets.Add above. Now update "Test Cases" rows 7 & 8.
# ---------------------------------------------------------------------------
$testCases = $wb.Worksheets.Item("Test Cases")

$testCases.Range("B7").Value2 = "OPQA-5781"
$testCases.Range("A7").Value2 = "PUBLONS008"
$testCases.Range("C7").Value2 = 'Verify that error message "Please enter your last name." whenever not enter any text in email field'
$testCases.Range("D7").Value2 = "Y"

$testCases.Range("A8").Value2 = "PUBLONS009"
$testCases.Range("B8").Value2 = "OPQA-5782&&OPQA-5783"
$testCases.Range("C8").Value2 = 'Verify that last name should be maximum of 50 characters long and these fields should not be empty.&&Verify that error message Last name is too long." whenever enter more than 50 characters'
$testCases.Range("D8").Value2 = "Y"

$testCases.Rows.Item(8).RowHeight = 30

$testCases.Activate() | Out-Null
$testCases.Range("C8").Select() | Out-Null
